$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114-159 down to 115-160.
$ws.Rows("114").Insert()

# Fill the newly inserted row 114 with the new record's data.
$ws.Range("A114").Value2 = 4
$ws.Range("B114").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value2 = "Los Lagos"
$ws.Range("D114").Value2 = 45093
$ws.Range("E114").Value2 = 10
$ws.Range("F114").Value2 = 100112022
$ws.Range("G114").Value2 = "Arveja Verde"
$ws.Range("H114").Value2 = "Perfection"
$ws.Range("I114").Value2 = "Primera"
$ws.Range("J114").Value2 = 40
$ws.Range("K114").Value2 = 42000
$ws.Range("L114").Value2 = 42000
$ws.Range("M114").Value2 = 42000
$ws.Range("N114").Value2 = "$/malla 25 kilos"
$ws.Range("O114").Value2 = "Provincia de Huasco"
$ws.Range("P114").Value2 = 1680
$ws.Range("Q114").Value2 = 25
$ws.Range("R114").Value2 = "Hortaliza"

# Ensure the date cell keeps the date/time number format used by the rest of column D.
$ws.Range("D114").NumberFormat = $ws.Range("D115").NumberFormat
